$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# Update the hours value in B21 (4 -> 9); dependent shared-formula cells in
# column C recalculate automatically.
$ws.Range("B21").Value = 9

# Update the active selection to match the saved view state.
$ws.Range("B22").Select()
